$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly data refresh inserts two brand-new daily records at the top of
# the data block (row 9/10), which pushes every existing record down by two
# rows (old row 9 -> new row 11, ..., old row 57 -> new row 59).
$ws.Rows("9:10").Insert()

# New record 1 (row 9)
$ws.Range("A9").Value = 6
$ws.Range("B9").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C9").Value = "Metropolitana"
$ws.Range("D9").Value = 45061
$ws.Range("E9").Value = 13
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100102
$ws.Range("H9").Value = "Cítricos"
$ws.Range("I9").Value = 100102006
$ws.Range("J9").Value = "Pomelo"
$ws.Range("K9").Value = "Start Ruby"
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 8
$ws.Range("N9").Value = 170000
$ws.Range("O9").Value = 170000
$ws.Range("P9").Value = 170000
$ws.Range("Q9").Value = "`$/bins (350 kilos)"
$ws.Range("R9").Value = "Región Metropolitana"
$ws.Range("S9").Value = 486
$ws.Range("T9").Value = 350

# New record 2 (row 10)
$ws.Range("A10").Value = 6
$ws.Range("B10").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C10").Value = "Metropolitana"
$ws.Range("D10").Value = 45061
$ws.Range("E10").Value = 13
$ws.Range("F10").Value = "Fruta"
$ws.Range("G10").Value = 100102
$ws.Range("H10").Value = "Cítricos"
$ws.Range("I10").Value = 100102006
$ws.Range("J10").Value = "Pomelo"
$ws.Range("K10").Value = "Start Ruby"
$ws.Range("L10").Value = "Segunda"
$ws.Range("M10").Value = 12
$ws.Range("N10").Value = 150000
$ws.Range("O10").Value = 150000
$ws.Range("P10").Value = 150000
$ws.Range("Q10").Value = "`$/bins (350 kilos)"
$ws.Range("R10").Value = "Región Metropolitana"
$ws.Range("S10").Value = 429
$ws.Range("T10").Value = 350
